# Day 2: Array, DP
# Applies:
#  - Sheet2 "Array": clean up 3-Num-Sum description, add "4 Num Sum" row,
#    add "Python Concepts" column (G)
#  - Sheet3 "Dynamic Programming": add "Min Number of Coins For Change" and
#    "Levenshtein Distance" rows
#  - New Sheet4 "Linked List" (header row only, copied styling from sheet3)
#  - View-state: selections per sheet + active sheet/tab back to Sheet1

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet2 : "Array"
# ---------------------------------------------------------------------------
$arr = $wb.Worksheets.Item("Array")

# Clean up the "3 Num Sum" description (remove embedded line-wraps from the
# original authoring, keep as a single paragraph).
$arr.Range("B3").Value = "Write a function that takes in a non-empty array of distinct integers and an integer representing a target sum. The function should find all triplets in the array that sum up to the target sum and return a two-dimensional array of all these triplets. The numbers in each triplet should be ordered in ascending order, and the triplets themselves should be ordered in ascending order with respect to the numbers they hold."
$arr.Rows.Item(3).RowHeight = 85.6

# New "Python Concepts" header column
$arr.Range("G1").Value = "Python Concepts"
$arr.Range("G1").Font.Bold = $true

# New "4 Num Sum" row of data (row 4 was blank before)
$arr.Range("A4").Value = "4 Num Sum "
$arr.Range("C4").Value = "[7, 6, 4, -1, 1, 2], target = 16 "
$arr.Range("D4").Value = "[[6, 7, 1, 2], [6, 7, -1, 4]]"
$arr.Range("E4").Value = "Average O(n^2) Worst O(n^3)"
$arr.Range("F4").Value = "O(n^2)"

$g4text = "extend ([6].extend([5] = [6, 5])),      append ([6].append([5] = [6, [5]]))"
$arr.Range("G4").Value = $g4text
$arr.Range("G4").Characters(1, 6).Font.Bold = $true
$arr.Range("G4").Characters(7, 34).Font.Bold = $false
$arr.Range("G4").Characters(41, 6).Font.Bold = $true
$arr.Range("G4").Characters(47, 29).Font.Bold = $false

$arr.Range("B4:G4").WrapText = $true
$arr.Rows.Item(4).RowHeight = 28.55

# Drop the now-superfluous extra blank spacer row (old row 6), keep row 5 as
# the single trailing blank spacer row.
$arr.Rows.Item(6).Delete()
$arr.Range("C5").WrapText = $true

# Column widths / new column
$arr.Columns.Item(5).ColumnWidth = 14.75
$arr.Columns.Item(5).WrapText = $true
$arr.Columns.Item(7).ColumnWidth = 35.125

# ---------------------------------------------------------------------------
# Sheet3 : "Dynamic Programming"
# ---------------------------------------------------------------------------
$dp = $wb.Worksheets.Item("Dynamic Programming")

$dp.Range("A4").Value = "Min Number of Coins For Change "
$dp.Range("B4").Value = "Given an array of positive integers representing coin denominations and a single non-negative integer representing a target amount of money, write a function that returns the smallest number of coins needed to make change for that target amount using the given coin denominations. If it's impossible to make change for the target amount, return -1"
$dp.Range("C4").Value = "n = 7, denoms = [1,5,10]"
$dp.Range("D4").Value = 3
$dp.Range("E4").Value = "O(nd)"
$dp.Range("F4").Value = "O(n)"
$dp.Range("B4:D4").WrapText = $true
$dp.Rows.Item(4).RowHeight = 114.15

$dp.Range("A5").Value = "Levenshtein Distance"
$dp.Range("B5").Value = "Write a function that takes in two strings and returns the minimum number of edit operations that need to be performed on the first string to obtain the second string."
$dp.Range("C5").Value = 'str1 = = "abc", str2= "yabd"   (insert "y"; substitute "c" for "d")'
$dp.Range("D5").Value = 2
$dp.Range("E5").Value = "O(nm)"
$dp.Range("F5").Value = "O(nm)"
$dp.Range("B5:F5").WrapText = $true
$dp.Rows.Item(5).RowHeight = 57.1

$dp.Columns.Item(4).ColumnWidth = 12.25
$dp.Columns.Item(4).WrapText = $true

# ---------------------------------------------------------------------------
# New Sheet4 : "Linked List" (copy DP's sheet for consistent base styling,
# then strip down to just the header row)
# ---------------------------------------------------------------------------
$dp.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ll = $wb.Worksheets.Item($wb.Worksheets.Count)
$ll.Name = "Linked List"
$ll.Rows.Item(2).Resize(2).Delete()

$ll.Range("A1").Value = "Problems"
$ll.Range("B1").Value = "Resources"
$ll.Range("C1").Value = "Input"
$ll.Range("D1").Value = "Output"
$ll.Range("E1").Value = "Time"
$ll.Range("F1").Value = "Space"
$ll.Range("G1").Value = "Python Concepts"
$ll.Range("A1:E1").WrapText = $true
$ll.Range("A1:G1").Font.Bold = $true

$ll.Columns.Item(1).ColumnWidth = 23.875
$ll.Columns.Item(1).WrapText = $true
$ll.Columns.Item(2).ColumnWidth = 39
$ll.Columns.Item(3).ColumnWidth = 30
$ll.Columns.Item(3).WrapText = $true
$ll.Columns.Item(4).ColumnWidth = 14.875
$ll.Columns.Item(4).WrapText = $true

# ---------------------------------------------------------------------------
# View state: per-sheet selection, restore Sheet1 as the active / selected
# tab (must be done last so it "sticks" as the saved active sheet).
# ---------------------------------------------------------------------------
$ll.Range("I5").Select()
$arr.Range("D11").Select()
$dp.Range("C11").Select()
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()
$sheet1.Range("C13").Select()

Write-Output "edit complete"
